$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.892.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.270.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.03%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.264.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("E10").Value = '  -5.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.571'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.71%  '
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '683.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.797.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.26%  '
$ws.Range("E16").Value = '  -3.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.043.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.268.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.94%  '
$ws.Range("E21").Value = '  -4.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.879'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.69%  '
$ws.Range("E26").Value = '  -4.50%  '
$ws.Range("E27").Value = '  -5.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '576.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.71%  '
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.791.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.26'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -15.75%  '
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -5.30%  '
$ws.Range("E42").Value = '  -4.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0651'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.68%  '
$ws.Range("E44").Value = '  -3.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.14%  '
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '127.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.27%  '

Write-Host "Updated cryptos list"
